$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C (IdxSG) values for rows 2..65 to the new numbering scheme ---
for ($row = 2; $row -le 65; $row++) {
    $val = 305020100 + 100 * ($row - 2)
    $ws.Cells.Item($row, 3).Value = $val
}

# --- Harmonise column B style for rows 3..65 (was s="3", becomes s="2", same as B2) ---
$ws.Range("B3:B65").VerticalAlignment = -4160

# --- Row 66 used to be a blank trailing row; it now becomes a real data row ---
$ws.Range("B66").Value = 2000000000
$ws.Range("B66").VerticalAlignment = -4160
$ws.Range("C66").Value = 305026500
$ws.Range("D66").Value = "ACTIVITES NON CLASSEES"

# --- Row 67 was a second blank trailing row; it is removed entirely ---
$ws.Rows("67").Delete()

# --- Selection moves from E59 (with frozen/scrolled top-left at A31) to a plain G7 selection ---
[void]$ws.Range("G7").Select()
